$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$lo = $ws.ListObjects.Item("Tabla1")

# --- Add the two new "Gap" columns to the table ---
$lo.ListColumns.Add() | Out-Null
$ws.Range("I2").Value2 = "Gap length (m)"
$lo.ListColumns.Add() | Out-Null
$ws.Range("J2").Value2 = "Gap Width (m)"

# --- Fill the data rows for the new columns ---
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = 0

# --- Column widths for the new columns (closest achievable to authored bestFit widths) ---
$ws.Columns.Item(9).ColumnWidth = 14.45
$ws.Columns.Item(10).ColumnWidth = 15.3

# --- Selection moves to the newly added header cell ---
$ws.Range("J2").Select() | Out-Null
